# Generate Report for Archive
#
# 1. Update the "Status" text from "Ready for handoff" to "In Translation"
#    everywhere it appears (Overview!E2:F2, zh-cn!C2, de-de!C2).
# 2. Narrow the "Status" column(s) (Overview E:F, zh-cn C, de-de C) to match
#    the refreshed report layout.

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"

# NOTE: this runtime's `.Value` *getter* doesn't resolve to the underlying
# scalar (it returns a reflection placeholder), so reads must go through
# `.Value2` (writes via `.Value` work fine). Also cast to [string] before
# comparing: some cells hold a real Boolean (e.g. "True"/"False" status
# flags) and PowerShell's `-eq` against a bare string operand coerces the
# string side to Boolean instead, producing false-positive matches.
foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    $rowCount = $used.Rows.Count
    $colCount = $used.Columns.Count
    for ($r = 1; $r -le $rowCount; $r++) {
        for ($c = 1; $c -le $colCount; $c++) {
            $cell = $used.Cells.Item($r, $c)
            $text = [string]$cell.Value2
            if ($text -eq $oldStatus) {
                $cell.Value = $newStatus
            }
        }
    }
}

# Resize the "Status" column on each sheet to the updated report width.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Columns.Item(3).ColumnWidth = 12.5

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
